$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed formatting for the new rows by copying from existing styled rows,
# then fill in the actual values (copy wipes the destination value too).
$fmtRow = 2
for ($r = 57; $r -le 68; $r++) {
  $ws.Range("A$fmtRow").Copy($ws.Range("A$r"))
  $ws.Range("B$fmtRow").Copy($ws.Range("B$r"))
  $ws.Range("C$fmtRow").Copy($ws.Range("C$r"))
  $ws.Range("D$fmtRow").Copy($ws.Range("D$r"))
  $ws.Range("E$fmtRow").Copy($ws.Range("E$r"))
  $ws.Range("F$fmtRow").Copy($ws.Range("F$r"))
}

$ws.Range("A57").Value = 'A098'
$ws.Range("B57").Value = '충청도_청주시_상당구'
$ws.Range("C57").Value = 'https://eminwon.cheongju.go.kr/emwp/jsp/ofr/OfrNotAncmtLSub.jsp?not_ancmt_se_code=01,04&yyyy=2020'
$ws.Range("D57").Value = '흥덕대교 보수보강공사 신기술·특정공법(표면보수)선정 기술제안서 제출안내 ...'
$ws.Range("E57").Value = 45652.0
$ws.Range("F57").Value = 45653.809641203705
$ws.Hyperlinks.Add($ws.Range("C57"), 'https://eminwon.cheongju.go.kr/emwp/jsp/ofr/OfrNotAncmtLSub.jsp?not_ancmt_se_code=01,04&yyyy=2020') | Out-Null
$ws.Range("E57:F57").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("A58").Value = 'A126'
$ws.Range("B58").Value = '전라도_전주시'
$ws.Range("C58").Value = 'https://eminwon.jeonju.go.kr/emwp/jsp/ofr/OfrNotAncmtLSub.jsp?not_ancmt_se_code=01,02,03,04&epcCheck=Y&recent_mm=60&list_gubun=A'
$ws.Range("D58").Value = '전주시 건축위원회 심의 기준 변경 공고'
$ws.Range("E58").Value = 45653.0
$ws.Range("F58").Value = 45653.809641203705
$ws.Hyperlinks.Add($ws.Range("C58"), 'https://eminwon.jeonju.go.kr/emwp/jsp/ofr/OfrNotAncmtLSub.jsp?not_ancmt_se_code=01,02,03,04&epcCheck=Y&recent_mm=60&list_gubun=A') | Out-Null
$ws.Range("E58:F58").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("A59").Value = 'A126'
$ws.Range("B59").Value = '전라도_전주시'
$ws.Range("C59").Value = 'https://eminwon.jeonju.go.kr/emwp/jsp/ofr/OfrNotAncmtLSub.jsp?not_ancmt_se_code=01,02,03,04&epcCheck=Y&recent_mm=60&list_gubun=A'
$ws.Range("D59").Value = '남부권 청소년센터 건립심의위원회 위원 공개모집 공고'
$ws.Range("E59").Value = 45653.0
$ws.Range("F59").Value = 45653.809641203705
$ws.Hyperlinks.Add($ws.Range("C59"), 'https://eminwon.jeonju.go.kr/emwp/jsp/ofr/OfrNotAncmtLSub.jsp?not_ancmt_se_code=01,02,03,04&epcCheck=Y&recent_mm=60&list_gubun=A') | Out-Null
$ws.Range("E59:F59").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("A60").Value = 'A126'
$ws.Range("B60").Value = '전라도_전주시'
$ws.Range("C60").Value = 'https://eminwon.jeonju.go.kr/emwp/jsp/ofr/OfrNotAncmtLSub.jsp?not_ancmt_se_code=01,02,03,04&epcCheck=Y&recent_mm=60&list_gubun=A'
$ws.Range("D60").Value = '전주시 용역과제심의위원회 위원 모집 공고'
$ws.Range("E60").Value = 45652.0
$ws.Range("F60").Value = 45653.809641203705
$ws.Hyperlinks.Add($ws.Range("C60"), 'https://eminwon.jeonju.go.kr/emwp/jsp/ofr/OfrNotAncmtLSub.jsp?not_ancmt_se_code=01,02,03,04&epcCheck=Y&recent_mm=60&list_gubun=A') | Out-Null
$ws.Range("E60:F60").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("A61").Value = 'A127'
$ws.Range("B61").Value = '전라도_정읍시'
$ws.Range("C61").Value = 'http://eminwon.jeongeup.go.kr/emwp/jsp/ofr/OfrNotAncmtL.jsp?not_ancmt_se_code=01,02,03,04,05,06,07'
$ws.Range("D61").Value = '매죽·매대 및 원종산 농촌마을하수도 설치사업  제안서 공법선정위원회 평가...'
$ws.Range("E61").Value = 45652.0
$ws.Range("F61").Value = 45653.809641203705
$ws.Hyperlinks.Add($ws.Range("C61"), 'http://eminwon.jeongeup.go.kr/emwp/jsp/ofr/OfrNotAncmtL.jsp?not_ancmt_se_code=01,02,03,04,05,06,07') | Out-Null
$ws.Range("E61:F61").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("A62").Value = 'A153'
$ws.Range("B62").Value = '전라도_장흥군'
$ws.Range("C62").Value = 'https://www.jangheung.go.kr/www/organization/news/notification'
$ws.Range("D62").Value = '문흥 빛의 거리 조성사업 제안서 평가 결과 공고새로운글'
$ws.Range("E62").Value = 45653.0
$ws.Range("F62").Value = 45653.809641203705
$ws.Hyperlinks.Add($ws.Range("C62"), 'https://www.jangheung.go.kr/www/organization/news/notification') | Out-Null
$ws.Range("E62:F62").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("A63").Value = 'A160'
$ws.Range("B63").Value = '경상도_구미시'
$ws.Range("C63").Value = 'https://www.gumi.go.kr/portal/saeol/gosi/list.do?seCode=01&mid=0401040000'
$ws.Range("D63").Value = '2025년 장애인 긴급돌보미 지원사업 수행기관 선정 심의 결과 공고'
$ws.Range("E63").Value = 45653.0
$ws.Range("F63").Value = 45653.809641203705
$ws.Hyperlinks.Add($ws.Range("C63"), 'https://www.gumi.go.kr/portal/saeol/gosi/list.do?seCode=01&mid=0401040000') | Out-Null
$ws.Range("E63:F63").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("A64").Value = 'A160'
$ws.Range("B64").Value = '경상도_구미시'
$ws.Range("C64").Value = 'https://www.gumi.go.kr/portal/saeol/gosi/list.do?seCode=01&mid=0401040000'
$ws.Range("D64").Value = '2025년 굿모닝 수요특강 위탁 용역 제안서 평가위원 명단 및 평가 결과 공개'
$ws.Range("E64").Value = 45652.0
$ws.Range("F64").Value = 45653.809641203705
$ws.Hyperlinks.Add($ws.Range("C64"), 'https://www.gumi.go.kr/portal/saeol/gosi/list.do?seCode=01&mid=0401040000') | Out-Null
$ws.Range("E64:F64").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("A65").Value = 'A177'
$ws.Range("B65").Value = '경상도_청도군'
$ws.Range("C65").Value = 'https://www.cheongdo.go.kr/portal/saeol/gosi/list.do?mid=0301020000&token=1703813401840'
$ws.Range("D65").Value = '청도군 청도읍 농촌중심지활성화사업 지역역량강화 용역 제안서 평가결과 공고'
$ws.Range("E65").Value = 45652.0
$ws.Range("F65").Value = 45653.809641203705
$ws.Hyperlinks.Add($ws.Range("C65"), 'https://www.cheongdo.go.kr/portal/saeol/gosi/list.do?mid=0301020000&token=1703813401840') | Out-Null
$ws.Range("E65:F65").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("A66").Value = 'A177'
$ws.Range("B66").Value = '경상도_청도군'
$ws.Range("C66").Value = 'https://www.cheongdo.go.kr/portal/saeol/gosi/list.do?mid=0301020000&token=1703813401840'
$ws.Range("D66").Value = '청도군 풍각면 기초생활거점조성사업 지역역량강화 용역 제안서 평가결과 공고'
$ws.Range("E66").Value = 45652.0
$ws.Range("F66").Value = 45653.809641203705
$ws.Hyperlinks.Add($ws.Range("C66"), 'https://www.cheongdo.go.kr/portal/saeol/gosi/list.do?mid=0301020000&token=1703813401840') | Out-Null
$ws.Range("E66:F66").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("A67").Value = 'A177'
$ws.Range("B67").Value = '경상도_청도군'
$ws.Range("C67").Value = 'https://www.cheongdo.go.kr/portal/saeol/gosi/list.do?mid=0301020000&token=1703813401840'
$ws.Range("D67").Value = '청도군 금천면 기초생활거점조성사업 지역역량강화 용역 제안서 평가결과 공고'
$ws.Range("E67").Value = 45652.0
$ws.Range("F67").Value = 45653.809641203705
$ws.Hyperlinks.Add($ws.Range("C67"), 'https://www.cheongdo.go.kr/portal/saeol/gosi/list.do?mid=0301020000&token=1703813401840') | Out-Null
$ws.Range("E67:F67").NumberFormat = "yyyy-mm-dd h:mm:ss"

$ws.Range("A68").Value = 'A177'
$ws.Range("B68").Value = '경상도_청도군'
$ws.Range("C68").Value = 'https://www.cheongdo.go.kr/portal/saeol/gosi/list.do?mid=0301020000&token=1703813401840'
$ws.Range("D68").Value = '『운문면 마일1리 농어촌취약지역생활여건개조사업 기본계획 수립 및 휴먼케어, 주민역량강화 용역』제안서 평가결과 공고'
$ws.Range("E68").Value = 45652.0
$ws.Range("F68").Value = 45653.809641203705
$ws.Hyperlinks.Add($ws.Range("C68"), 'https://www.cheongdo.go.kr/portal/saeol/gosi/list.do?mid=0301020000&token=1703813401840') | Out-Null
$ws.Range("E68:F68").NumberFormat = "yyyy-mm-dd h:mm:ss"
